# Apply the "Card21" sheet update: add a new service-log entry (row 32)
# and backfill previously-blank placeholder cells with "nan" text, as
# produced by the source export pipeline.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card21")

# --- Row 29: previously blank placeholder cells become literal "nan" text ---
$ws.Range("L29").Value = "nan"
$ws.Range("M29").Value = "nan"
$ws.Range("N29").Value = "nan"
$ws.Range("O29").Value = "nan"

# --- Row 31: previously blank placeholder cell becomes literal "nan" text ---
$ws.Range("O31").Value = "nan"

# --- New Row 32: new event entry for Card21 ---
# Column A holds the card number as text ("21"), matching the rest of the sheet.
$ws.Range("A32").NumberFormat = "@"
$ws.Range("A32").Value = "21"
$ws.Range("A32").ClearFormats()

# Columns B-K stay blank (no measurement data supplied for this event) but are
# still present as cells on the row.
$ws.Range("B32:K32").NumberFormat = "General"

# Columns L-O carry the new maintenance event details.
$ws.Range("L32").Value = "13/12/2025"
$ws.Range("M32").Value = "تكسير سنون فرشة الفلاتس البطئ"
$ws.Range("N32").Value = "تم تغيير فرشة الفلاتس البطئ عند طن (967)"
$ws.Range("O32").Value = "حسام السيد/سعيد فرج"
